$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. ANTI-COX II row 10: H10 stock 0:5 -> 1:5
$ws.Range("H10").Value = "1:5"

# 2. EREC row 32: H32, P32, Q32
$ws.Range("H32").Value = "3:12"
$ws.Range("P32").Value = "23.0400"
$ws.Range("Q32").Value = "0:2"

# 3. RICHI PANTHENOL row 56: H56
$ws.Range("H56").Value = "1:0"

# 4. STRINGAZOLE row 60: H60
$ws.Range("H60").Value = "1:0"

# 5. URIVIN-N row 66: H66
$ws.Range("H66").Value = "6:0"

Write-Output "done phase 1"
